$d = $word.ActiveDocument

# NOTE ON ORDERING
# ------------------------------------------------------------------
# Several blocks of text move from one place in the document to
# another (e.g. the "Fomentar a cultura..." paragraph moves into the
# Docente(s) list, the "Norma de recuperação:" value moves into the
# "Programa resumido" paragraph, etc). Word's Find/Replace
# (Replace = wdReplaceAll) matches *every* occurrence of the search
# text in the document, so the operations below are ordered such that
# a given source text is always located and consumed (via Find) while
# it is still unique, *before* any later step re-introduces an equal
# copy of that same text elsewhere. The one paragraph that gains new
# runs (Docente(s) list) is rebuilt last, using literal text (no
# Find), so it cannot clash with earlier/later searches.
# ------------------------------------------------------------------

$wdReplaceOne = 1

# 1) "Ativação: 01/01/2024" -> "Ativação: Semestral"
$d.Content.Find.Execute("Ativação: 01/01/2024", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: Semestral", $wdReplaceOne) | Out-Null

# 2) "Norma de recuperação:" value (currently "NF = (MF + PR)/2...") is replaced
#    with the new extensionist-activity description. Do this before step 6
#    (which reintroduces the "NF = (MF + PR)/2..." text elsewhere).
$d.Content.Find.Execute("NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação", $true, $false, $false, $false, $false, $true, 1, $false, "Esta atividade é denominada Engenharia e Negócios – Oficina de Empreendedorismo e Inovação.A atividade consiste na realização de uma oficina de Empreendedorismo e Inovação com estudantes do Ensino Médio. Tal oficina poderá ocorrer em escolas de ensino médio ou em organizações sociais ou representativas das comunidades da cidade de Lorena (ex: associações de bairros).Etapas:1.Planejamento da Oficina: definição dos temas (na área de empreendedorismo e inovação) a serem desenvolvidos, que pode incluir concursos de ideias, desafios de negócios, exposição de conteúdos, entre outras atividades, com a consequente preparação dos materiais (slides, vídeos, exercícios etc.) que serão utilizados nas oficinas. Os alunos serão os responsáveis por agendar a realização das oficinas com as escolas ou com outras organizações previamente aprovadas pelo professor da disciplina.2.Realização da Oficina: poderá ser aplicada em uma ou duas etapas (dias diferentes), somando no mínimo 4 horas totais de aplicação.3.Preparação de Relatos em Vídeo: criação de um vídeo relatando o desenvolvimento da oficina, com os aprendizados adquiridos, para ser disponibilizado para a comunidade.4.Autoavaliação pelo Grupo: avaliar os resultados da avaliação da atividade aplicada aos estudantes do ensino médio, para identificar o aprendizado e os pontos a melhorar para as próximas oficinas.", $wdReplaceOne) | Out-Null

# 3) "Bibliografia" body paragraph (currently the "BLANK, Steve Gary..." reference
#    list) is replaced with the survey-method text. Do this before step 7 (which
#    reintroduces the "BLANK, Steve Gary..." text elsewhere).
$d.Content.Find.Execute("BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.DOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. FILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. HISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.", $true, $false, $false, $false, $false, $true, 1, $false, "Será realizada uma pesquisa de satisfação com os participantes da oficina e para o responsável pela atividade na escola de ensino médio ou organização. Após a pesquisa, o grupo de estudantes da disciplina, fará uma análise dos resultados e uma autoavaliação e discutirá tais resultados com o professor da disciplina", $wdReplaceOne) | Out-Null

# 4) "Avaliação" list paragraph: "Método:" value text. Do this before step 9
#    (Docente list) reintroduces the same "Aulas expositivas..." text.
$d.Content.Find.Execute("Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", $true, $false, $false, $false, $false, $true, 1, $false, "O grupo social alvo são estudantes do ensino médio em escolas públicas e/ou nas comunidades (associações de bairros etc.) da cidade de Lorena/SP.", $wdReplaceOne) | Out-Null

# 5) "Avaliação" list paragraph: "Critério:" value text. Do this before step 9
#    (Docente list) reintroduces the same "Média Aritmética..." text.
$d.Content.Find.Execute("Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas", $true, $false, $false, $false, $false, $true, 1, $false, "São objetivos da atividade Extensionista:- Disseminar a cultura empreendedora e inovadora aos estudantes de Ensino Médio;- Estimular os estudantes de Ensino Médio para o desenvolvimento de sua capacidade empreendedora, a busca de oportunidades, a geração do autoemprego e o desenvolvimento de atitudes empreendedoras e criativas.", $wdReplaceOne) | Out-Null

# 6) "Programa resumido" body paragraph: replace with "NF = (MF + PR)/2, ...".
#    (The only earlier occurrence of this text was already consumed in step 2.)
$d.Content.Find.Execute("Características do Comportamento Empreendedor; Modelo de Negócios; Produto mínimo viável; Plano de Negócios.", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação", $wdReplaceOne) | Out-Null

# 7) "Programa" body paragraph: replace the long syllabus text with the
#    bibliography list. (The only earlier occurrence of this text was already
#    consumed in step 3.)
$d.Content.Find.Execute("1. Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2. Estratégia, Inovação e Marketing.3. Design Thinking.4. Modelo de Negócios (Business Model Canvas e Lean Startup - Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.5. Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente. Prototipação rápida.6. Gestão de processos e Gerenciamento ágil de projetos.7. Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.8. Proposta da criação de uma startup, do modelo de negócios ao plano de negócios, incluindo a montagem do produto mínimo viável e uma rodada de PITCH.9. Desenvolvimento de atividade prática extensionista (produção de conteúdo digital sobre empreendedorismo e inovação)10. Visita (viagem didática complementar) a um ambiente de inovação e empreendedorismo (ex. incubadora/aceleradora ou parque tecnológico), para compreender o desenvolvimento dos processos de empreendedorismo e inovação.", $true, $false, $false, $false, $false, $true, 1, $false, "BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.DOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. FILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. HISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.", $wdReplaceOne) | Out-Null

# 8) "Objetivos" body paragraph: the long "Fomentar a cultura..." text is
#    replaced by "01/01/2025". Do this before step 9 reintroduces the same
#    "Fomentar a cultura..." text elsewhere.
$d.Content.Find.Execute("Fomentar a cultura do empreendedorismo e da Inovação; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup com uma proposta de produto/serviço inovador ao longo do semestre.", $true, $false, $false, $false, $false, $true, 1, $false, "01/01/2025", $wdReplaceOne) | Out-Null

# 9) "Docente(s) Responsável(eis)" list paragraph: currently a single run
#    "11079086 - Herlandí de Souza Andrade". Expand it into six runs separated
#    by manual line breaks, pulling in text that moved from the paragraphs
#    above. This uses literal strings (InsertBefore/InsertAfter, not Find), so
#    it cannot collide with the Find-based replacements performed earlier.
$brk = [char]11
$paras = $d.Paragraphs
$docentePara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
  $p = $paras.Item($i)
  if ($p.Range.Text.StartsWith("11079086 - Herland")) {
    $docentePara = $p
    break
  }
}
$rng = $docentePara.Range
# Insert the "Fomentar..." text (+break) before the existing text.
$rng.InsertBefore("Fomentar a cultura do empreendedorismo e da Inovação; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup com uma proposta de produto/serviço inovador ao longo do semestre." + $brk)

# Re-find the paragraph (its range changed) and append the remaining new runs
# after the existing "11079086 - ..." run, each separated by a manual line
# break, with no trailing break after the very last one.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
  $p = $paras.Item($i)
  if ($p.Range.Text.IndexOf("11079086 - Herland") -ge 0) {
    $docentePara = $p
    break
  }
}
$endRng = $docentePara.Range
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1)
$endRng.InsertAfter($brk + "Características do Comportamento Empreendedor; Modelo de Negócios; Produto mínimo viável; Plano de Negócios." + $brk + "1. Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2. Estratégia, Inovação e Marketing.3. Design Thinking.4. Modelo de Negócios (Business Model Canvas e Lean Startup - Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.5. Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente. Prototipação rápida.6. Gestão de processos e Gerenciamento ágil de projetos.7. Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.8. Proposta da criação de uma startup, do modelo de negócios ao plano de negócios, incluindo a montagem do produto mínimo viável e uma rodada de PITCH.9. Desenvolvimento de atividade prática extensionista (produção de conteúdo digital sobre empreendedorismo e inovação)10. Visita (viagem didática complementar) a um ambiente de inovação e empreendedorismo (ex. incubadora/aceleradora ou parque tecnológico), para compreender o desenvolvimento dos processos de empreendedorismo e inovação." + $brk + "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras." + $brk + "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas")
